$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the cells that were removed from the "notification" rows
# (D = image, G = link in the header row) for rows 2 and 3.
$ws.Range("D2").ClearContents()
$ws.Range("G2").ClearContents()
$ws.Range("D3").ClearContents()
$ws.Range("G3").ClearContents()

# Selection moved from G6 to G3
$null = $ws.Range("G3").Select()
